$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as text (matches source data which
# stores prices as plain strings, e.g. "26.864.27", "1.001"), so Excel does not
# reinterpret them as numbers/dates and strip formatting such as trailing zeros.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.864.27"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.863.43"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "304.92"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5047"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").Value = "0.07163"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "0.8920"
$ws.Range("E10").Value = "  +0.79%  "
$ws.Range("D11").Value = "20.61"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "0.07514"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("D13").Value = "1.861.27"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "95.06"
$ws.Range("E14").Value = "  +6.85%  "
$ws.Range("D15").Value = "5.226"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("D17").Value = "0.000008510"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "14.21"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "0.9994"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "26.921.01"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "5.025"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "2.094.78"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "10.36"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").Value = "6.405"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "147.95"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "1.778"
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").Value = "17.88"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "2.071"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "113.19"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("D30").Value = "4.689"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").Value = "4.660"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "0.09175"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "0.05138"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").Value = "0.7493"
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("D35").Value = "2.976"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").Value = "1.153"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "3.193"
$ws.Range("E37").Value = "  +5.23%  "
$ws.Range("D38").Value = "2.562"
$ws.Range("E38").Value = "  +4.48%  "
$ws.Range("D39").Value = "0.01995"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").Value = "0.5578"
$ws.Range("E40").Value = "  +5.88%  "
$ws.Range("D41").Value = "1.070"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "6.583"
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").Value = "115.88"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").Value = "8.539"
$ws.Range("E44").Value = "  +3.80%  "
$ws.Range("D45").Value = "0.1470"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "0.4692"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "0.9995"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "10.08"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "1.555"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "36.72"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "62.99"
$ws.Range("E51").Value = "  -1.30%  "
